$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 78; existing rows 78-92 shift down to 80-94.
$ws.Rows("78:79").Insert()

# --- Row 78: new weekly price record (2022-03-07, Maracuyá, Especial) ---
$ws.Cells.Item(78, 1).Value = 1
$ws.Cells.Item(78, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(78, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(78, 4).Value = 44627
$ws.Cells.Item(78, 5).Value = 15
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100108
$ws.Cells.Item(78, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(78, 9).Value = 100108003
$ws.Cells.Item(78, 10).Value = "Maracuyá"
$ws.Cells.Item(78, 11).Value = "Sin especificar"
$ws.Cells.Item(78, 12).Value = "Especial"
$ws.Cells.Item(78, 13).Value = 120
$ws.Cells.Item(78, 14).Value = 21000
$ws.Cells.Item(78, 15).Value = 22000
$ws.Cells.Item(78, 16).Value = 21500
$ws.Cells.Item(78, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(78, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(78, 19).Value = 1075
$ws.Cells.Item(78, 20).Value = 20

# --- Row 79: new weekly price record (2022-03-07, Maracuyá, Primera) ---
$ws.Cells.Item(79, 1).Value = 1
$ws.Cells.Item(79, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(79, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(79, 4).Value = 44627
$ws.Cells.Item(79, 5).Value = 15
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100108
$ws.Cells.Item(79, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(79, 9).Value = 100108003
$ws.Cells.Item(79, 10).Value = "Maracuyá"
$ws.Cells.Item(79, 11).Value = "Sin especificar"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 130
$ws.Cells.Item(79, 14).Value = 19000
$ws.Cells.Item(79, 15).Value = 20000
$ws.Cells.Item(79, 16).Value = 19500
$ws.Cells.Item(79, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(79, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 19).Value = 975
$ws.Cells.Item(79, 20).Value = 20
